$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted above the old row 335, pushing the
# existing rows 335-361 down to 336-362 (dimension grows from A1:R361 to
# A1:R362). Insert a blank row at 335 so everything below shifts down.
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new data point.
$ws.Range("A335").Value = 8
$ws.Range("B335").Value = "Terminal La Palmera de La Serena"
$ws.Range("C335").Value = "Coquimbo"
$ws.Range("D335").Value = 45166
$ws.Range("E335").Value = 4
$ws.Range("F335").Value = 100112037
$ws.Range("G335").Value = "Cebollín"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 1200
$ws.Range("K335").Value = 1000
$ws.Range("L335").Value = 1200
$ws.Range("M335").Value = 1100
$ws.Range("N335").Value = "$/paquete 6 unidades"
$ws.Range("O335").Value = "Provincia del Elquí"
$ws.Range("P335").Value = 183
$ws.Range("Q335").Value = 6
$ws.Range("R335").Value = "Hortaliza"

# Give the new date cell the same date/time number format used by the rest
# of column D (style index 2 in the original workbook).
$ws.Range("D335").NumberFormat = $ws.Range("D336").NumberFormat
